$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# Row 2 (W8)
$wsForecast.Range("H2").Value = 1.67
$wsForecast.Range("L2").Value = 1.18

# Row 3 (W9)
$wsForecast.Range("H3").Value = 0.51
$wsForecast.Range("I3").Value = "Low"
$wsForecast.Range("L3").Value = 1

# Row 4 (W10)
$wsForecast.Range("L4").Value = 0.99

# Row 5 (W11)
$wsForecast.Range("L5").Value = 1.15

# Row 6 (W12)
$wsForecast.Range("L6").Value = 0.8100000000000001

# Row 7 (W13)
$wsForecast.Range("L7").Value = 1.09

# Row 8 (W14)
$wsForecast.Range("L8").Value = 0.85

# Row 9 (W15)
$wsForecast.Range("L9").Value = 1.17

# Row 11 (W17)
$wsForecast.Range("L11").Value = 1.1

# Row 12 (W18)
$wsForecast.Range("L12").Value = 0.93

# Row 13 (W19)
$wsForecast.Range("L13").Value = 0.82

# Row 14 (W20)
$wsForecast.Range("L14").Value = 0.99

# Row 15 (W21)
$wsForecast.Range("L15").Value = 1.11

# Row 16 (W22)
$wsForecast.Range("L16").Value = 1

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.13

# Summary sheet: Total Forecast (16 Weeks)
$wsSummary.Range("B9").Value = "10"
